$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 61, shifting existing rows 61:153 down to 62:154
$ws.Rows("61").Insert()

# Populate the newly inserted row 61 with the new weekly data point
$ws.Range("A61").Value = 5
$ws.Range("B61").Value = "Macroferia Regional de Talca"
$ws.Range("C61").Value = "Maule"
$ws.Range("D61").Value = 44665
$ws.Range("E61").Value = 7
$ws.Range("F61").Value = 100112017
$ws.Range("G61").Value = "Apio"
$ws.Range("H61").Value = "Americana (o)"
$ws.Range("I61").Value = "Primera"
$ws.Range("J61").Value = 500
$ws.Range("K61").Value = 8000
$ws.Range("L61").Value = 8000
$ws.Range("M61").Value = 8000
$ws.Range("N61").Value = "$/docena de matas"
$ws.Range("O61").Value = "Provincia del Elquí"
$ws.Range("P61").Value = 1333
$ws.Range("Q61").Value = 6
$ws.Range("R61").Value = "Hortaliza"
